$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Convert Excel serial date 45221 to a real date so the cell keeps its
# existing date formatting/style while only the stored value changes.
$newDate = [DateTime]::FromOADate(45221)

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
